$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "11 - Gaspésie – Îles-de-la-Madeleine" case count (row 13)
$ws.Range("B13").Value = 19138

# Update "Total" case count (row 22)
$ws.Range("B22").Value = 372476
